$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the price/volume columns to be treated as plain text so that
# numeric-looking values (e.g. "1.000", "243.84") are not reinterpreted
# by Excel as numbers/dates. We restore the original ("Normal") cell
# style afterwards so no stray number formats remain applied.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "29.939.86"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "1.892.67"
$ws.Range("E3").Value = "  -0.05%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "0.7740"
$ws.Range("E5").Value = "  -2.72%  "
$ws.Range("D6").Value = "243.84"
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "0.3139"
$ws.Range("E8").Value = "  -0.86%  "
$ws.Range("D9").Value = "25.75"
$ws.Range("E9").Value = "  +1.01%  "
$ws.Range("D10").Value = "0.07358"
$ws.Range("E10").Value = "  +4.35%  "
$ws.Range("D11").Value = "0.08064"
$ws.Range("E11").Value = "  -0.20%  "
$ws.Range("D12").Value = "0.7722"
$ws.Range("E12").Value = "  +0.33%  "
$ws.Range("D13").Value = "5.504"
$ws.Range("E13").Value = "  +2.74%  "
$ws.Range("D14").Value = "1.847.12"
$ws.Range("E14").Value = "  -2.64%  "
$ws.Range("D15").Value = "94.14"
$ws.Range("E15").Value = "  +1.62%  "
$ws.Range("E16").Value = "  +3.70%  "
$ws.Range("D17").Value = "29.917.19"
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("D18").Value = "14.01"
$ws.Range("E18").Value = "  +0.95%  "
$ws.Range("D19").Value = "246.63"
$ws.Range("E19").Value = "  +0.79%  "
$ws.Range("E20").Value = "  +1.90%  "
$ws.Range("D21").Value = "8.154"
$ws.Range("E21").Value = "  -2.55%  "
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("B23").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C23").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D23").Value = "2.136.77"
$ws.Range("E23").Value = "  -0.80%  "
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").Value = "0.1574"
$ws.Range("E25").Value = "  -3.93%  "
$ws.Range("D26").Value = "9.445"
$ws.Range("E26").Value = "  +0.95%  "
$ws.Range("D27").Value = "161.97"
$ws.Range("E27").Value = "  -2.38%  "
$ws.Range("D28").Value = "18.76"
$ws.Range("E28").Value = "  +0.29%  "
$ws.Range("D29").Value = "2.026"
$ws.Range("E29").Value = "  -1.61%  "
$ws.Range("E30").Value = "  +1.75%  "
$ws.Range("D31").Value = "1.542"
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("D32").Value = "4.470"
$ws.Range("E32").Value = "  +0.53%  "
$ws.Range("D33").Value = "0.05565"
$ws.Range("E33").Value = "  -2.21%  "
$ws.Range("D34").Value = "4.064"
$ws.Range("E34").Value = "  +0.45%  "
$ws.Range("E35").Value = "  -1.89%  "
$ws.Range("E36").Value = "  +1.91%  "
$ws.Range("D37").Value = "1.002"
$ws.Range("E37").Value = "  +0.20%  "
$ws.Range("D38").Value = "2.683"
$ws.Range("E38").Value = "  +2.35%  "
$ws.Range("D39").Value = "0.01930"
$ws.Range("E39").Value = "  +1.02%  "
$ws.Range("D40").Value = "2.790"
$ws.Range("E40").Value = "  +0.15%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "0.4476"
$ws.Range("E41").Value = "  +1.55%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "74.35"
$ws.Range("E42").Value = "  +2.59%  "
$ws.Range("D43").Value = "1.101.09"
$ws.Range("E43").Value = "  +6.60%  "
$ws.Range("D44").Value = "6.021"
$ws.Range("E44").Value = "  +3.42%  "
$ws.Range("E45").Value = "  +1.18%  "
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("E47").Value = "  +0.83%  "
$ws.Range("D48").Value = "102.47"
$ws.Range("E48").Value = "  -0.57%  "
$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").Value = "7.544"
$ws.Range("E49").Value = "  +1.55%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "9.772"
$ws.Range("E50").Value = "  -2.20%  "
$ws.Range("D51").Value = "2.999"
$ws.Range("E51").Value = "  +3.19%  "

# Restore default styling on the touched range (clears the temporary
# text-number-format applied above).
$dataRange.Style = "Normal"

